$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1) updates to column F (想去人数)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 438
$ws1.Range("F7").Value = 568
$ws1.Range("F9").Value = 6845
$ws1.Range("F11").Value = 100
$ws1.Range("F16").Value = 16285
$ws1.Range("F17").Value = 1603
$ws1.Range("F19").Value = 334
$ws1.Range("F20").Value = 188
$ws1.Range("F22").Value = 11407
$ws1.Range("F24").Value = 1056
$ws1.Range("F25").Value = 4494
$ws1.Range("F30").Value = 320

# Sheet "全部类型" (sheet4) updates to column F (想去人数)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 438
$ws4.Range("F7").Value = 568
$ws4.Range("F10").Value = 6845
$ws4.Range("F12").Value = 100
$ws4.Range("F18").Value = 16285
$ws4.Range("F19").Value = 1603
$ws4.Range("F21").Value = 334
$ws4.Range("F22").Value = 188
$ws4.Range("F26").Value = 11407
$ws4.Range("F28").Value = 1056
$ws4.Range("F29").Value = 4494
$ws4.Range("F34").Value = 320

$wb.Save()
